$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that need to switch from the "general" style (s=7) to the
# "date" style (s=8, numFmtId 14) because they now hold date values.
# Copy the number-format/style from an existing date-formatted cell
# (G3, style s="8") and paste only the formatting onto each target cell.
$cellsNeedingDateStyle = @("H8","G9","H9","G10","H10","G11","H11","G12","H12")
foreach ($cellAddr in $cellsNeedingDateStyle) {
    $ws.Range("G3").Copy()
    $ws.Range($cellAddr).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Row 8
$ws.Range("G8").Value = 43400
$ws.Range("H8").Value = 43401

# Row 9
$ws.Range("G9").Value = 43401
$ws.Range("H9").Value = 43401

# Row 10
$ws.Range("F10").Value = 43406
$ws.Range("G10").Value = "29/10/2018"
$ws.Range("H10").Value = 43403

# Row 11
$ws.Range("E11").Value = 43407
$ws.Range("F11").Value = 43408
$ws.Range("G11").Value = 43405
$ws.Range("H11").Value = 43406

# Row 12
$ws.Range("G12").Value = 43407
$ws.Range("H12").Value = 43407

# Keep selection consistent with the authored workbook state.
[void]$ws.Range("H14").Select()
